# Daily update of the Valais COVID-19 figures workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One additional positive case recorded retro-actively (row 296).
$ws.Range("C296").Value = 81

# New daily figures for the most recent days (rows 474-476).
$ws.Range("C474").Value = 5
$ws.Range("C475").Value = 5
$ws.Range("C476").Value = 1

# L476/M476 are formatted as Text ("@"); a plain .Value assignment would
# store the number as a text string. Temporarily switch the format to
# General so the underlying value is written as a genuine number, then
# restore the original Text format (this reuses the pre-existing style).
$deathsRange = $ws.Range("L476:M476")
$deathsRange.NumberFormat = "General"
$ws.Range("L476").Value = 0
$ws.Range("M476").Value = 0
$deathsRange.NumberFormat = "@"

# Move the active selection of the frozen (bottom-right) pane to A2.
$ws.Activate()
$ws.Range("A2").Select()
